$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Hardcoded input values that changed ---
$ws.Range("B8").Value = 123
$ws.Range("B12").Value = 337
$ws.Range("B13").Value = 9999
$ws.Range("B47").Value = 760

# --- Unit label text changes ---
$ws.Range("C18").Value = "WattHour / GHash"
$ws.Range("C19").Value = "kW / Watt"

# --- Column B width change (12 -> 16.5 raw xml units) ---
$ws.Columns.Item(2).ColumnWidth = 15.666666666666666

# --- View/selection state ---
$ws.Range("E26").Select()
$excel.ActiveWindow.ScrollRow = 15
$excel.ActiveWindow.ScrollColumn = 1
